$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 45076
$ws.Range("M2").Value = 150
$ws.Range("O2").Value = 11000
$ws.Range("P2").Value = 10467
$ws.Range("R2").Value = 'Provincia de Curicó'
$ws.Range("S2").Value = 582

# Row 3
$ws.Range("D3").Value = 45037
$ws.Range("M3").Value = 250
$ws.Range("N3").Value = 9000
$ws.Range("O3").Value = 9500
$ws.Range("P3").Value = 9200
$ws.Range("Q3").Value = '$/caja 18 kilos granel'
$ws.Range("R3").Value = 'Provincia de Curicó'
$ws.Range("S3").Value = 511

# Row 4
$ws.Range("D4").Value = 45027
$ws.Range("N4").Value = 9000
$ws.Range("O4").Value = 10000
$ws.Range("P4").Value = 9500
$ws.Range("S4").Value = 528

# Row 5
$ws.Range("D5").Value = 45029
$ws.Range("M5").Value = 100
$ws.Range("O5").Value = 10000
$ws.Range("P5").Value = 9500
$ws.Range("Q5").Value = '$/bandeja 18 kilos granel'
$ws.Range("R5").Value = 'Región de O''Higgins'
$ws.Range("S5").Value = 528

# Row 6
$ws.Range("D6").Value = 45092
$ws.Range("M6").Value = 110
$ws.Range("N6").Value = 10000
$ws.Range("O6").Value = 11000
$ws.Range("P6").Value = 10455
$ws.Range("Q6").Value = '$/bandeja 18 kilos granel'
$ws.Range("R6").Value = 'Provincia de Curicó'
$ws.Range("S6").Value = 581

# Row 7
$ws.Range("D7").Value = 45126
$ws.Range("N7").Value = 14000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 14500
$ws.Range("Q7").Value = '$/bandeja 18 kilos granel'
$ws.Range("S7").Value = 806
$ws.Range("T7").Value = 18

# Row 8
$ws.Range("D8").Value = 44999
$ws.Range("N8").Value = 12000
$ws.Range("O8").Value = 12000
$ws.Range("P8").Value = 12000
$ws.Range("S8").Value = 667

# Row 9
$ws.Range("L9").Value = 'Segunda'
$ws.Range("N9").Value = 10000
$ws.Range("O9").Value = 10000
$ws.Range("P9").Value = 10000
$ws.Range("S9").Value = 556

# Row 10
$ws.Range("D10").Value = 44299
$ws.Range("L10").Value = 'Primera'
$ws.Range("O10").Value = 11000
$ws.Range("P10").Value = 10500
$ws.Range("Q10").Value = '$/caja 18 kilos granel'
$ws.Range("R10").Value = 'Región del Maule'
$ws.Range("S10").Value = 583

# Row 11
$ws.Range("D11").Value = 44299
$ws.Range("L11").Value = 'Segunda'
$ws.Range("M11").Value = 50
$ws.Range("O11").Value = 9000
$ws.Range("P11").Value = 9000
$ws.Range("R11").Value = 'Región del Maule'
$ws.Range("S11").Value = 500

# Row 12
$ws.Range("D12").Value = 44698
$ws.Range("Q12").Value = '$/caja 18 kilos granel'

# Row 13
$ws.Range("D13").Value = 44307
$ws.Range("L13").Value = 'Primera'
$ws.Range("N13").Value = 10000
$ws.Range("O13").Value = 10000
$ws.Range("P13").Value = 10000
$ws.Range("S13").Value = 556

# Row 14
$ws.Range("D14").Value = 44307
$ws.Range("L14").Value = 'Segunda'
$ws.Range("N14").Value = 8000
$ws.Range("O14").Value = 8000
$ws.Range("P14").Value = 8000
$ws.Range("S14").Value = 444

# Row 15
$ws.Range("D15").Value = 45050
$ws.Range("M15").Value = 140
$ws.Range("N15").Value = 11000
$ws.Range("O15").Value = 12000
$ws.Range("P15").Value = 11429
$ws.Range("Q15").Value = '$/caja 18 kilos empedrada'
$ws.Range("R15").Value = 'Región de O''Higgins'
$ws.Range("S15").Value = 635

# Row 16
$ws.Range("D16").Value = 45107
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = 11000
$ws.Range("P16").Value = 11000
$ws.Range("Q16").Value = '$/caja 18 kilos empedrada'
$ws.Range("R16").Value = 'Región del Maule'
$ws.Range("S16").Value = 611

# Row 17
$ws.Range("D17").Value = 44316
$ws.Range("N17").Value = 9000
$ws.Range("O17").Value = 10000
$ws.Range("P17").Value = 9500
$ws.Range("Q17").Value = '$/caja 18 kilos granel'
$ws.Range("S17").Value = 528

# Row 18
$ws.Range("D18").Value = 45079
$ws.Range("M18").Value = 270
$ws.Range("N18").Value = 11000
$ws.Range("O18").Value = 12000
$ws.Range("P18").Value = 11444
$ws.Range("Q18").Value = '$/caja 18 kilos granel'
$ws.Range("S18").Value = 636

# Row 19
$ws.Range("D19").Value = 45034
$ws.Range("M19").Value = 220
$ws.Range("N19").Value = 8500
$ws.Range("O19").Value = 9000
$ws.Range("P19").Value = 8727
$ws.Range("Q19").Value = '$/caja 18 kilos granel'
$ws.Range("R19").Value = 'Región de O''Higgins'
$ws.Range("S19").Value = 485

# Row 20
$ws.Range("D20").Value = 45041
$ws.Range("N20").Value = 11000
$ws.Range("O20").Value = 12000
$ws.Range("P20").Value = 11500
$ws.Range("S20").Value = 639

# Row 21
$ws.Range("D21").Value = 44358
$ws.Range("M21").Value = 100
$ws.Range("P21").Value = 11500
$ws.Range("Q21").Value = '$/caja 18 kilos granel'
$ws.Range("S21").Value = 639

# Row 22
$ws.Range("D22").Value = 44776
$ws.Range("M22").Value = 50
$ws.Range("O22").Value = 10000
$ws.Range("P22").Value = 10000
$ws.Range("Q22").Value = '$/bandeja 18 kilos granel'
$ws.Range("R22").Value = 'Región de O''Higgins'
$ws.Range("S22").Value = 556

# Row 23
$ws.Range("D23").Value = 44776
$ws.Range("N23").Value = 8000
$ws.Range("O23").Value = 8000
$ws.Range("P23").Value = 8000
$ws.Range("Q23").Value = '$/bandeja 18 kilos granel'
$ws.Range("R23").Value = 'Región de O''Higgins'
$ws.Range("S23").Value = 444

# Row 24
$ws.Range("D24").Value = 45013

# Row 25
$ws.Range("D25").Value = 44363
$ws.Range("M25").Value = 100
$ws.Range("N25").Value = 9000
$ws.Range("P25").Value = 9500
$ws.Range("Q25").Value = '$/caja 15 kilos empedrada'
$ws.Range("S25").Value = 633
$ws.Range("T25").Value = 15

# Row 26
$ws.Range("D26").Value = 45128
$ws.Range("L26").Value = 'Primera'
$ws.Range("N26").Value = 12000
$ws.Range("O26").Value = 12000
$ws.Range("P26").Value = 12000
$ws.Range("S26").Value = 667

# Row 27
$ws.Range("D27").Value = 44272
$ws.Range("M27").Value = 100
$ws.Range("N27").Value = 9000
$ws.Range("O27").Value = 10000
$ws.Range("P27").Value = 9500
$ws.Range("Q27").Value = '$/caja 15 kilos granel'
$ws.Range("S27").Value = 633
$ws.Range("T27").Value = 15

# Row 28
$ws.Range("L28").Value = 'Segunda'
$ws.Range("M28").Value = 50
$ws.Range("N28").Value = 8000
$ws.Range("O28").Value = 8000
$ws.Range("P28").Value = 8000
$ws.Range("S28").Value = 533

# Row 29
$ws.Range("D29").Value = 44425
$ws.Range("L29").Value = 'Primera'
$ws.Range("M29").Value = 100
$ws.Range("N29").Value = 12000
$ws.Range("O29").Value = 13000
$ws.Range("P29").Value = 12500
$ws.Range("Q29").Value = '$/bandeja 18 kilos granel'
$ws.Range("S29").Value = 694
$ws.Range("T29").Value = 18

# Row 30
$ws.Range("D30").Value = 45014
$ws.Range("N30").Value = 9000
$ws.Range("O30").Value = 10000
$ws.Range("P30").Value = 9500
$ws.Range("Q30").Value = '$/bandeja 18 kilos granel'
$ws.Range("S30").Value = 528
